$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order timestamps regenerated) ---
$wb.Worksheets.Item("GNG_TO-16511687425650241").Name = "GNG_TO-16512555704071147"
$wb.Worksheets.Item("NB_TO-1651168745246005").Name = "NB_TO-16512555734150155"
$wb.Worksheets.Item("RS_TO-1651168745246434").Name = "RS_TO-16512555734170165"
$wb.Worksheets.Item("TOL_TO-16511687453087118").Name = "TOL_TO-16512555734640176"
$wb.Worksheets.Item("vSAT_TO-16511687453850935").Name = "vSAT_TO-16512555735407321"

# --- Sheet 1: GNG ---
$ws = $wb.Worksheets.Item("GNG_TO-16512555704071147")
$ws.Range("B2").Value = "go_stims-16512555703713923.csv"
$ws.Range("B3").Value = "GNG_stims-16512555703910575.csv"
$ws.Range("B4").Value = "go_stims-16512555703930588.csv"
$ws.Range("B5").Value = "GNG_stims-16512555704061122.csv"

# --- Sheet 2: NB ---
$ws = $wb.Worksheets.Item("NB_TO-16512555734150155")
$ws.Range("B2").Value = "ZB-match_8-1651255571103362.csv"
$ws.Range("B3").Value = "ZB-match_7-16512555708811545.csv"
$ws.Range("B4").Value = "TB-1651255572169777.csv"
$ws.Range("B5").Value = "TB-16512555720162358.csv"
$ws.Range("B6").Value = "OB-16512555716165235.csv"
$ws.Range("B7").Value = "OB-16512555716600084.csv"
$ws.Range("B8").Value = "ZB-match_5-1651255570813676.csv"
$ws.Range("B9").Value = "OB-16512555717837315.csv"
$ws.Range("B10").Value = "TB-1651255573394015.csv"

# --- Sheet 3: RS ---
$ws = $wb.Worksheets.Item("RS_TO-16512555734170165")
$ws.Range("B2").Value = "eyes closed"
$ws.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws = $wb.Worksheets.Item("TOL_TO-16512555734640176")
$ws.Range("B2").Value = "MM_stims-16512555734310172.csv"
$ws.Range("B3").Value = "ZM_stims-1651255573419018.csv"
$ws.Range("B4").Value = "MM_stims-16512555734470181.csv"
$ws.Range("B5").Value = "ZM_stims-16512555734320226.csv"
$ws.Range("B6").Value = "MM_stims-16512555734620202.csv"
$ws.Range("B7").Value = "ZM_stims-16512555734480178.csv"

# --- Sheet 5: vSAT ---
$ws = $wb.Worksheets.Item("vSAT_TO-16512555735407321")
$ws.Range("B2").Value = "SAT_stims-16512555734680192.csv"
$ws.Range("B3").Value = "vSAT_stims-1651255573510026.csv"
$ws.Range("B4").Value = "SAT_stims-16512555734950252.csv"
$ws.Range("B5").Value = "vSAT_stims-1651255573525705.csv"
